$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assign as literal text, preserving the original General/default style.
# (a plain Range.Value = "7.71" would be auto-coerced to a Double by Excel
# type inference, losing formatting like trailing zeros or the "." digit
# grouping used elsewhere in this sheet - so we force Text first, then
# restore the Normal/General style so no stray style index is left behind).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "63.605.05"
Set-TextValue "E2" "  +2.69%  "
Set-TextValue "D3" "3.477.32"
Set-TextValue "E3" "  +1.59%  "
Set-TextValue "D5" "581.91"
Set-TextValue "E5" "  +0.58%  "
Set-TextValue "D6" "147.68"
Set-TextValue "E6" "  +1.50%  "
Set-TextValue "D7" "3.478.01"
Set-TextValue "E7" "  +1.60%  "
Set-TextValue "E8" "  -0.16%  "
Set-TextValue "E9" "  +0.61%  "
Set-TextValue "D10" "7.71"
Set-TextValue "E10" "  +0.90%  "
Set-TextValue "E11" "  +1.59%  "
Set-TextValue "E12" "  +4.83%  "
Set-TextValue "D13" "4.071.85"
Set-TextValue "E13" "  +1.63%  "
Set-TextValue "E14" "  +5.27%  "
Set-TextValue "E15" "  +2.42%  "
Set-TextValue "D16" "3.479.02"
Set-TextValue "E16" "  +1.75%  "
Set-TextValue "D17" "0.0000173"
Set-TextValue "E17" "  +1.45%  "
Set-TextValue "D18" "63.503.45"
Set-TextValue "E18" "  +2.59%  "
Set-TextValue "D19" "6.36"
Set-TextValue "E19" "  +2.79%  "
Set-TextValue "D20" "14.44"
Set-TextValue "E20" "  +3.84%  "
Set-TextValue "D21" "9.36"
Set-TextValue "E21" "  +1.66%  "
Set-TextValue "D22" "390.60"
Set-TextValue "E22" "  -0.03%  "
Set-TextValue "D23" "0.566"
Set-TextValue "E23" "  +2.37%  "
Set-TextValue "D24" "75.10"
Set-TextValue "E24" "  +0.89%  "
Set-TextValue "E25" "  -0.22%  "
Set-TextValue "D26" "3.621.42"
Set-TextValue "E26" "  +1.75%  "
Set-TextValue "E27" "  +0.93%  "
Set-TextValue "D28" "0.181"
Set-TextValue "E28" "  -4.95%  "
Set-TextValue "D29" "7.64"
Set-TextValue "E29" "  +2.18%  "
Set-TextValue "E30" "  +0.36%  "
Set-TextValue "D31" "8.26"
Set-TextValue "E31" "  +2.88%  "
Set-TextValue "E32" "  -0.11%  "
Set-TextValue "E34" "  -4.11%  "
Set-TextValue "D35" "23.58"
Set-TextValue "E35" "  +0.18%  "
Set-TextValue "D36" "5.33"
Set-TextValue "E36" "  +1.15%  "
Set-TextValue "D37" "7.15"
Set-TextValue "E37" "  +2.37%  "
Set-TextValue "E38" "  +8.25%  "
Set-TextValue "D39" "31.64"
Set-TextValue "E39" "  +9.63%  "
Set-TextValue "D40" "169.85"
Set-TextValue "E40" "  +0.46%  "
Set-TextValue "D41" "3.514.52"
Set-TextValue "E41" "  +1.68%  "
Set-TextValue "D42" "0.0766"
Set-TextValue "E42" "  +1.37%  "
Set-TextValue "D43" "0.801"
Set-TextValue "E43" "  +1.70%  "
Set-TextValue "D44" "1.74"
Set-TextValue "E44" "  +3.64%  "
Set-TextValue "D45" "42.44"
Set-TextValue "E45" "  -0.80%  "
Set-TextValue "D46" "1.22"
Set-TextValue "E46" "  +3.42%  "
Set-TextValue "D47" "4.43"
Set-TextValue "E47" "  -0.54%  "
Set-TextValue "D48" "2.611.21"
Set-TextValue "E48" "  +3.30%  "
Set-TextValue "E49" "  +9.28%  "
Set-TextValue "D50" "23.16"
Set-TextValue "E50" "  +1.29%  "
Set-TextValue "D51" "6.79"
